$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '91.483.51'
$ws.Range("E2").Value = '  -0.30%  '
$ws.Range("D3").Value = '3.109.79'
$ws.Range("E3").Value = '  +0.01%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '243.16'
$ws.Range("E5").Value = '  -0.55%  '
$ws.Range("D6").Value = '615.84'
$ws.Range("E6").Value = '  -0.80%  '
$ws.Range("E7").Value = '  -2.95%  '
$ws.Range("D8").Value = '0.384'
$ws.Range("E8").Value = '  +3.43%  '
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  -0.14%  '
$ws.Range("D10").Value = '3.106.63'
$ws.Range("E10").Value = '  +13.87%  '
$ws.Range("D11").Value = '0.740'
$ws.Range("E11").Value = '  -2.44%  '
$ws.Range("E12").Value = '  +0.44%  '
$ws.Range("D13").Value = '0.0000249'
$ws.Range("E13").Value = '  -1.59%  '
$ws.Range("D14").Value = '5.60'
$ws.Range("E14").Value = '  +2.40%  '
$ws.Range("D15").Value = '34.46'
$ws.Range("E15").Value = '  -2.80%  '
$ws.Range("D16").Value = '91.284.50'
$ws.Range("E16").Value = '  -0.39%  '
$ws.Range("D18").Value = '3.108.94'
$ws.Range("E18").Value = '  -0.20%  '
$ws.Range("D19").Value = '3.70'
$ws.Range("E19").Value = '  -0.85%  '
$ws.Range("E20").Value = '  +1.00%  '
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '9.28'
$ws.Range("E22").Value = '  +1.84%  '
$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").Value = '443.87'
$ws.Range("E23").Value = '  -0.73%  '
$ws.Range("E24").Value = '  -7.52%  '
$ws.Range("D25").Value = '5.61'
$ws.Range("E25").Value = '  -0.19%  '
$ws.Range("D26").Value = '87.99'
$ws.Range("E26").Value = '  -2.76%  '
$ws.Range("D27").Value = '11.64'
$ws.Range("E27").Value = '  -2.54%  '
$ws.Range("B29").Value = 'Hedera'
$ws.Range("C29").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D29").Value = '0.142'
$ws.Range("E29").Value = '  +27.90%  '
$ws.Range("B30").Value = 'Dai'
$ws.Range("C30").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D30").Value = '0.996'
$ws.Range("E30").Value = '  -0.49%  '
$ws.Range("D31").Value = '0.230'
$ws.Range("E31").Value = '  -2.00%  '
$ws.Range("D32").Value = '0.166'
$ws.Range("E32").Value = '  -9.45%  '
$ws.Range("E33").Value = '  +3.94%  '
$ws.Range("E34").Value = '  -1.53%  '
$ws.Range("E35").Value = '  -1.84%  '
$ws.Range("D36").Value = '7.62'
$ws.Range("E36").Value = '  -0.58%  '
$ws.Range("D37").Value = '26.12'
$ws.Range("E37").Value = '  -1.55%  '
$ws.Range("D38").Value = '4.06'
$ws.Range("E38").Value = '  -1.33%  '
$ws.Range("E39").Value = '  +1.10%  '
$ws.Range("D40").Value = '487.38'
$ws.Range("E40").Value = '  -0.79%  '
$ws.Range("D41").Value = '1.30'
$ws.Range("E41").Value = '  -0.83%  '
$ws.Range("E42").Value = '  +3.32%  '
$ws.Range("D43").Value = '3.40'
$ws.Range("E43").Value = '  -6.49%  '
$ws.Range("D44").Value = '22.19'
$ws.Range("E44").Value = '  +0.10%  '
$ws.Range("E45").Value = '  -0.02%  '
$ws.Range("D46").Value = '159.22'
$ws.Range("E46").Value = '  +2.95%  '
$ws.Range("E47").Value = '  +0.60%  '
$ws.Range("D48").Value = '1.89'
$ws.Range("E48").Value = '  -1.06%  '
$ws.Range("E49").Value = '  -0.63%  '
$ws.Range("D50").Value = '44.05'
$ws.Range("E50").Value = '  -1.20%  '
$ws.Range("D51").Value = '4.37'
$ws.Range("E51").Value = '  -3.03%  '
